# Correction type pour génération à partir fsh ea4a6f04ed193a83290686b2f69a3f9cd2e7f4ad
#
# - Metadata sheet: "Name" row (B4) gets the value "PaysnaissanceVs"
# - Metadata sheet: "Date" row (B8) value is refreshed to the new timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B4").Value = "PaysnaissanceVs"
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
